$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains text formatting for numeric-looking values,
# matching the source data which stores all Price values as text.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '22.451.28'
$ws.Range("E2").Value = '  +0.17%  '
$ws.Range("D3").Value = '1.573.09'
$ws.Range("E3").Value = '  +0.63%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("E5").Value = '  +0.04%  '
$ws.Range("D6").Value = '288.03'
$ws.Range("E6").Value = '  +0.89%  '
$ws.Range("D7").Value = '0.3699'
$ws.Range("E7").Value = '  +1.72%  '
$ws.Range("B8").Value = 'OKB'
$ws.Range("C8").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D8").Value = '47.31'
$ws.Range("E8").Value = '  -2.44%  '
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").Value = '0.3322'
$ws.Range("E9").Value = '  -0.49%  '
$ws.Range("D10").Value = '1.155'
$ws.Range("E10").Value = '  +2.43%  '
$ws.Range("D11").Value = '0.07507'
$ws.Range("E11").Value = '  +1.25%  '
$ws.Range("E12").Value = '  +0.07%  '
$ws.Range("E13").Value = '  -0.14%  '
$ws.Range("D14").Value = '5.936'
$ws.Range("E14").Value = '  +0.09%  '
$ws.Range("E15").Value = '  +0.41%  '
$ws.Range("D16").Value = '1.562.71'
$ws.Range("E16").Value = '  -0.09%  '
$ws.Range("D17").Value = '0.00001115'
$ws.Range("E17").Value = '  +0.86%  '
$ws.Range("D18").Value = '88.44'
$ws.Range("E18").Value = '  +0.23%  '
$ws.Range("D19").Value = '0.06724'
$ws.Range("E19").Value = '  +0.64%  '
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Value = '6.391'
$ws.Range("E20").Value = '  +0.42%  '
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").Value = '0.9997'
$ws.Range("E21").Value = '  -0.06%  '
$ws.Range("E22").Value = '  +2.32%  '
$ws.Range("E23").Value = '  +0.02%  '
$ws.Range("D24").Value = '22.433.57'
$ws.Range("E24").Value = '  +0.09%  '
$ws.Range("D25").Value = '2.385'
$ws.Range("E25").Value = '  -1.12%  '
$ws.Range("D26").Value = '2.631'
$ws.Range("E26").Value = '  +2.43%  '
$ws.Range("D27").Value = '150.86'
$ws.Range("E27").Value = '  +0.68%  '
$ws.Range("D28").Value = '19.59'
$ws.Range("E28").Value = '  +0.86%  '
$ws.Range("D29").Value = '4.959'
$ws.Range("E29").Value = '  -0.72%  '
$ws.Range("D30").Value = '125.04'
$ws.Range("E30").Value = '  +1.60%  '
$ws.Range("D31").Value = '1.740.90'
$ws.Range("E31").Value = '  +0.07%  '
$ws.Range("D32").Value = '1.099'
$ws.Range("E32").Value = '  +3.12%  '
$ws.Range("D33").Value = '6.082'
$ws.Range("E33").Value = '  -1.06%  '
$ws.Range("D34").Value = '1.990'
$ws.Range("E34").Value = '  -0.16%  '
$ws.Range("D35").Value = '9.903'
$ws.Range("E35").Value = '  +0.99%  '
$ws.Range("D36").Value = '0.08328'
$ws.Range("E36").Value = '  +0.96%  '
$ws.Range("E37").Value = '  +1.91%  '
$ws.Range("D38").Value = '1.307'
$ws.Range("E38").Value = '  +0.03%  '
$ws.Range("D39").Value = '0.06378'
$ws.Range("E39").Value = '  -0.12%  '
$ws.Range("D40").Value = '0.2217'
$ws.Range("E40").Value = '  +0.34%  '
$ws.Range("D41").Value = '5.333'
$ws.Range("E41").Value = '  -0.06%  '
$ws.Range("E42").Value = '  +1.88%  '
$ws.Range("D43").Value = '0.6241'
$ws.Range("E43").Value = '  +2.52%  '
$ws.Range("B44").Value = 'Frax'
$ws.Range("C44").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D44").Value = '0.9999'
$ws.Range("E44").Value = '  +0.01%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '13.99'
$ws.Range("E45").Value = '  +1.33%  '
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = '0.6062'
$ws.Range("E46").Value = '  +5.20%  '
$ws.Range("B47").Value = 'PancakeSwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D47").Value = '3.774'
$ws.Range("E47").Value = '  +0.42%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = '2.043'
$ws.Range("E48").Value = '  +1.52%  '
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").Value = '124.86'
$ws.Range("E49").Value = '  -0.03%  '
$ws.Range("B50").Value = 'EOS'
$ws.Range("C50").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D50").Value = '1.206'
$ws.Range("E50").Value = '  -0.76%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '0.07195'
$ws.Range("E51").Value = '  -0.25%  '
